$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the grand-total formula in F2: was SUM(D2:D9), now SUM(D:D) ---
$ws.Range("F2").Formula = "=SUM(D:D)"

# --- Extend the shared "=Bn*Cn" formula pattern down through the new rows ---
for ($r = 8; $r -le 17; $r++) {
    $ws.Range("D$r").Formula = "=B$r*C$r"
}

# --- Row 7 keeps its own content; rows 8-10 get new part descriptions/values ---

# Give C12 (new "Roller Bearing" row) its no-decimals currency format FIRST so
# that it claims the lower new style index, matching the new-style ordering.
$ws.Range("C12").NumberFormat = "\$#,##0_);[Red](\$#,##0)"

$ws.Range("A8").Value = "`n11 GA. (.120 thick)`nHot Rolled Steel Sheet (2'X4')"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 77
$ws.Range("A8").Font.Name = "Arial"
$ws.Range("A8").Font.Size = 11
$ws.Range("A8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 43.5

$ws.Range("A9").Value = "5/8"" Hot Rolled A-36 Steel Round"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 5.14
$ws.Range("A9").Font.Name = "Arial"
$ws.Range("A9").Font.Size = 11
$ws.Range("A9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 29.25

$ws.Range("A10").Value = "1/2"" Hot Rolled A-36 Steel Round"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 6.6
$ws.Range("A10").Font.Color = 0
$ws.Range("A10").Font.Name = "Arial"
$ws.Range("A10").Font.Size = 11

# --- Brand new rows 11-18: additional parts ---
$ws.Range("A11").Value = "Bearing"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = 11

$ws.Range("A12").Value = "Roller Bearing"
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 10

$ws.Range("A13").Value = "Pulley"
$ws.Range("B13").Value = 4
$ws.Range("C13").Value = 16

$ws.Range("A14").Value = "Driving Belt"
$ws.Range("B14").Value = 2
$ws.Range("C14").Value = 6

$ws.Range("A15").Value = "Gears"
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = 22

$ws.Range("A16").Value = "Gears"
$ws.Range("B16").Value = 2
$ws.Range("C16").Value = 13

$ws.Range("A17").Value = "5"" wheels/hubs"
$ws.Range("B17").Value = 4
$ws.Range("C17").Value = 16

$ws.Range("A18").Value = "Screws"
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 20
$ws.Range("D18").Formula = "=B18*C18"

# --- Formatting: numbers in column C/D use the 2-decimal currency format ---
$ws.Range("C2:C11").NumberFormat = "\$#,##0.00_);[Red](\$#,##0.00)"
$ws.Range("C13:C18").NumberFormat = "\$#,##0.00_);[Red](\$#,##0.00)"
$ws.Range("D2:D18").NumberFormat = "\$#,##0.00_);[Red](\$#,##0.00)"
$ws.Range("F2").NumberFormat = "\$#,##0.00_);[Red](\$#,##0.00)"

# --- View state: scrolled down a bit, selection moved to the grand-total cell ---
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("F2").Select()

$wb.Save()
